# Atualização de bases das ligas, do dia: 17-02-2024 às 22:47
#
# This edit:
#  1) Swaps the full data (columns B:AC) of rows 112 and 113 (id/A stays put).
#  2) Swaps the full data (columns B:AC) of rows 263 and 265 (id/A stays put).
#  3) Shifts the data (columns B:AA) of rows 296->295, 297->296, 298->297
#     (id/A stays put in each destination row), then removes the now
#     duplicated trailing row 298 entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($ws, $row1, $row2, $firstCol, $lastCol) {
    $addr1 = $firstCol + $row1 + ":" + $lastCol + $row1
    $addr2 = $firstCol + $row2 + ":" + $lastCol + $row2
    $rng1 = $ws.Range($addr1)
    $rng2 = $ws.Range($addr2)
    $count = $rng1.Columns.Count
    for ($i = 1; $i -le $count; $i++) {
        $c1 = $rng1.Cells.Item(1, $i)
        $c2 = $rng2.Cells.Item(1, $i)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

function Copy-RowData($ws, $srcRow, $dstRow, $firstCol, $lastCol) {
    $addrSrc = $firstCol + $srcRow + ":" + $lastCol + $srcRow
    $addrDst = $firstCol + $dstRow + ":" + $lastCol + $dstRow
    $src = $ws.Range($addrSrc)
    $dst = $ws.Range($addrDst)
    $count = $src.Columns.Count
    for ($i = 1; $i -le $count; $i++) {
        $dst.Cells.Item(1, $i).Value2 = $src.Cells.Item(1, $i).Value2
    }
}

# 1) Rows 112 & 113 swap (columns B..AC)
Swap-RowData $ws 112 113 "B" "AC"

# 2) Rows 263 & 265 swap (columns B..AC)
Swap-RowData $ws 263 265 "B" "AC"

# 3) Rows 295..298 shift up by one (columns B..AA), then drop row 298
Copy-RowData $ws 296 295 "B" "AA"
Copy-RowData $ws 297 296 "B" "AA"
Copy-RowData $ws 298 297 "B" "AA"

$ws.Rows(298).Delete()
